# Quest/dialogue sheet header renewal: capitalize header row (name -> Name,
# index -> Index, dialogueIndex -> DialogueIndex, needItem -> NeedItem,
# needCount -> NeedCount, rewardItem -> RewardItem, rewardCount -> RewardCount).
# Data rows are left untouched. Also move the active selection to G1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Index"
$ws.Range("C1").Value = "DialogueIndex"
$ws.Range("D1").Value = "NeedItem"
$ws.Range("E1").Value = "NeedCount"
$ws.Range("F1").Value = "RewardItem"
$ws.Range("G1").Value = "RewardCount"

# Move selection to G1 (matches the saved cursor position in the diff).
$ws.Range("G1").Select()
